$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple value updates (no reordering) ---
# Row 4: Estados Unidos
$ws.Range("B4").Value = 1375401
$ws.Range("C4").Value = 7763
$ws.Range("E4").Value = 1036223
$ws.Range("G4").Value = 354
$ws.Range("H4").Value = 81141

# Row 12: Turquia
$ws.Range("B12").Value = 139771
$ws.Range("C12").Value = 1114
$ws.Range("D12").Value = 95780
$ws.Range("E12").Value = 40150
$ws.Range("F12").Value = 1126
$ws.Range("G12").Value = 55
$ws.Range("H12").Value = 3841

# Row 32: Emiratos Arabes Unidos
$ws.Range("B32").Value = 18878
$ws.Range("C32").Value = 680
$ws.Range("D32").Value = 5381
$ws.Range("E32").Value = 13296
$ws.Range("G32").Value = 3
$ws.Range("H32").Value = 201

# --- Reorder rows 113-115: Mali moves above Crucero, ahead of Uruguay ---
# Row 113 becomes Mali (with updated stats)
$ws.Range("A113").Value = "Mali"
$ws.Range("B113").Value = 712
$ws.Range("C113").Value = 8
$ws.Range("D113").Value = 377
$ws.Range("E113").Value = 296
$ws.Range("F113").Value = 0
$ws.Range("G113").Value = 1
$ws.Range("H113").Value = 39

# Row 114 becomes Crucero (old Crucero stats)
$ws.Range("A114").Value = "Crucero"
$ws.Range("B114").Value = 712
$ws.Range("C114").Value = 0
$ws.Range("D114").Value = 651
$ws.Range("E114").Value = 48
$ws.Range("F114").Value = 4
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = 13

# Row 115 becomes Uruguay (old Uruguay stats)
$ws.Range("A115").Value = "Uruguay"
$ws.Range("B115").Value = 707
$ws.Range("C115").Value = 0
$ws.Range("D115").Value = 517
$ws.Range("E115").Value = 171
$ws.Range("F115").Value = 8
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = 19

# --- Reorder rows 192-193: Nueva Caledonia moves above Belice ---
# Row 192 becomes Nueva Caledonia
$ws.Range("A192").Value = "Nueva Caledonia"
$ws.Range("D192").Value = 18
$ws.Range("H192").Value = 0

# Row 193 becomes Belice
$ws.Range("A193").Value = "Belice"
$ws.Range("D193").Value = 16
$ws.Range("H193").Value = 2
